$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column B ("Age"), shifting old B..H to C..I ---
$ws.Columns("B:B").Insert()

# --- 2. Header row ---
$ws.Range("B1").Value = "Age"
$ws.Range("J1").Value = "Free_Text"

# --- 3. Age values for each data row (col B) ---
$ages = @{
  2 = 35; 3 = 26; 4 = 27; 5 = 34; 6 = 30; 7 = 32; 8 = 32; 9 = 25; 10 = 25;
  11 = 24; 12 = 25; 13 = 31; 14 = 27; 15 = 31; 16 = 29; 17 = 34; 18 = 32;
  19 = 26; 20 = 20; 21 = 24
}
foreach ($r in $ages.Keys) {
  $ws.Cells.Item($r, 2).Value = $ages[$r]
}

# --- 4. Previously-missing previous_experience values (now col G) ---
$ws.Range("G4").Value = "{'Mitz Marak': 1}"
$ws.Range("G15").Value = "{'Cafe Yehoshua': 3}"
$ws.Range("G17").Value = "{'Hagadir': 3}"
$ws.Range("G18").Value = "{'Mitz Marak': 4}"

# --- 5. Column widths ---
$ws.Range("A:B").ColumnWidth = 12.77734375
$ws.Range("C:C").ColumnWidth = 17.88671875
$ws.Range("D:D").ColumnWidth = 12.21875
$ws.Range("E:E").ColumnWidth = 13.21875
$ws.Range("F:F").ColumnWidth = 20.21875
$ws.Range("G:G").ColumnWidth = 22.21875
$ws.Range("I:I").ColumnWidth = 11.6640625
$ws.Range("J:J").ColumnWidth = 21.44140625

# --- 6. Sheet view: scroll + selection ---
$ws.Range("J2").Select()

# --- 7. Hyperlink now lives on C2 (was B2, shifted by the insert) ---

Write-Host "edit complete"
